# cross-refs, captions, other formatting
#
# The "Heading 1" style picks up a forced page break before each
# occurrence (so every top-level heading starts its own page).
$d = $word.ActiveDocument

$heading1 = $d.Styles("Heading1")
$heading1.ParagraphFormat.PageBreakBefore = $true
